{"js": "// Office.js (Word JavaScript API) edit script\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// ---------------------------------------------------------------\n// 1) Collapse the three \"CORE COMPETENCIES\" detail paragraphs into a\n//    single condensed paragraph.\n// ---------------------------------------------------------------\nconst productMarketingCoreText =\n  \"Product Marketing Core: Market Intelligence & Competitive Analysis \\u2022 Product Positioning & Messaging Development \\u2022 Go-to-Market Strategy & Product Launch Management \\u2022 Customer Segmentation & Buyer Persona Development\";\nconst researchAnalyticsText =\n  \"Research & Analytics: Survey Methodology & Customer Insights \\u2022 Market Research Design & Implementation \\u2022 Competitive Intelligence & SWOT Analysis \\u2022 A/B Testing & Conversion Optimization\";\nconst communicationTechText =\n  \"Communication & Technology: Strategic Messaging & Narrative Development \\u2022 Technical Concept Translation for Business Audiences \\u2022 Data Visualization & Reporting (Tableau, PowerBI, d3.js) \\u2022 Client Relationship Management & Business Development\";\n\nlet coreIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === productMarketingCoreText) {\n    coreIdx = i;\n    break;\n  }\n}\n\nif (coreIdx === -1) {\n  throw new Error(\"Could not find 'Product Marketing Core' paragraph\");\n}\n\n// Confirm the next two paragraphs are the ones we expect before removing.\nif (\n  items[coreIdx + 1].text !== researchAnalyticsText ||\n  items[coreIdx + 2].text !== communicationTechText\n) {\n  throw new Error(\"Unexpected paragraphs following 'Product Marketing Core'\");\n}\n\nitems[coreIdx].insertText(\n  \"Product Marketing Core \\u2022 Research & Analytics \\u2022 Communication & Technology\",\n  Word.InsertLocation.replace\n);\nitems[coreIdx + 1].delete();\nitems[coreIdx + 2].delete();\n\nawait context.sync();\n\n// ---------------------------------------------------------------\n// 2) Insert a new \"TECHNICAL SKILLS\" section (heading + 3 detail\n//    paragraphs) right after the \"Managed national polling team...\"\n//    bullet, which is the last bullet before the closing sentence.\n// ---------------------------------------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"\\u2022 Managed national polling team of five data analysts for consumer insights and market intelligence\";\n\nlet anchorIdx = -1;\nconst refreshedItems = paragraphs.items;\nfor (let i = 0; i < refreshedItems.length; i++) {\n  if (refreshedItems[i].text === anchorText) {\n    anchorIdx = i;\n    break;\n  }\n}\n\nif (anchorIdx === -1) {\n  throw new Error(\"Could not find anchor paragraph 'Managed national polling team...'\");\n}\n\nconst anchorParagraph = refreshedItems[anchorIdx];\n\n// Insert paragraphs in reverse order, each time right after the anchor,\n// so the final order reading top-to-bottom is: heading, then the 3\n// detail lines.\nconst commTechPara = anchorParagraph.insertParagraph(\n  \"COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Client Relationship Management & Business Development\",\n  Word.InsertLocation.after\n);\nconst researchPara = anchorParagraph.insertParagraph(\n  \"RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; A/B Testing & Conversion Optimization\",\n  Word.InsertLocation.after\n);\nconst productPara = anchorParagraph.insertParagraph(\n  \"PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development\",\n  Word.InsertLocation.after\n);\nconst headingPara = anchorParagraph.insertParagraph(\n  \"TECHNICAL SKILLS\",\n  Word.InsertLocation.after\n);\nheadingPara.style = \"Heading 2\";\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n\n$productMarketingCoreText = \"Product Marketing Core: Market Intelligence & Competitive Analysis \" + $bullet + \" Product Positioning & Messaging Development \" + $bullet + \" Go-to-Market Strategy & Product Launch Management \" + $bullet + \" Customer Segmentation & Buyer Persona Development\"\n$researchAnalyticsText = \"Research & Analytics: Survey Methodology & Customer Insights \" + $bullet + \" Market Research Design & Implementation \" + $bullet + \" Competitive Intelligence & SWOT Analysis \" + $bullet + \" A/B Testing & Conversion Optimization\"\n$communicationTechText = \"Communication & Technology: Strategic Messaging & Narrative Development \" + $bullet + \" Technical Concept Translation for Business Audiences \" + $bullet + \" Data Visualization & Reporting (Tableau, PowerBI, d3.js) \" + $bullet + \" Client Relationship Management & Business Development\"\n\n# ---------------------------------------------------------------\n# 1) Collapse the three \"CORE COMPETENCIES\" detail paragraphs into a\n#    single condensed paragraph.\n# ---------------------------------------------------------------\n$coreIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -eq ($productMarketingCoreText + \"`r\")) {\n        $coreIndex = $i\n        break\n    }\n}\n\nif ($coreIndex -eq -1) {\n    throw \"Could not find 'Product Marketing Core' paragraph\"\n}\n\nif ($d.Paragraphs($coreIndex + 1).Range.Text -ne ($researchAnalyticsText + \"`r\")) {\n    throw \"Unexpected paragraph following 'Product Marketing Core' (Research & Analytics)\"\n}\nif ($d.Paragraphs($coreIndex + 2).Range.Text -ne ($communicationTechText + \"`r\")) {\n    throw \"Unexpected paragraph following 'Product Marketing Core' (Communication & Technology)\"\n}\n\n# NOTE: Range.Text assignments intentionally omit the trailing \"`r\" --\n# including it would insert an *additional* empty paragraph instead of\n# just replacing this paragraph's content in place.\n$d.Paragraphs($coreIndex).Range.Text = \"Product Marketing Core \" + $bullet + \" Research & Analytics \" + $bullet + \" Communication & Technology\"\n$d.Paragraphs($coreIndex + 1).Range.Delete()\n$d.Paragraphs($coreIndex + 1).Range.Delete()\n\n# ---------------------------------------------------------------\n# 2) Insert a new \"TECHNICAL SKILLS\" section (heading + 3 detail\n#    paragraphs) right after the \"Managed national polling team...\"\n#    bullet, which is the last bullet before the closing sentence.\n# ---------------------------------------------------------------\n$anchorText = $bullet + \" Managed national polling team of five data analysts for consumer insights and market intelligence\"\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -eq ($anchorText + \"`r\")) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph 'Managed national polling team...'\"\n}\n\n# The paragraph immediately after the anchor is where we insert the\n# new content (i.e. insert 4 new paragraphs before it, one at a time).\n$insertBeforeIndex = $anchorIndex + 1\n\n$d.Paragraphs($insertBeforeIndex).Range.InsertParagraphBefore()\n$d.Paragraphs($insertBeforeIndex).Range.Text = \"TECHNICAL SKILLS\"\n$d.Paragraphs($insertBeforeIndex).Style = \"Heading 2\"\n$insertBeforeIndex++\n\n$d.Paragraphs($insertBeforeIndex).Range.InsertParagraphBefore()\n$d.Paragraphs($insertBeforeIndex).Range.Text = \"PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development\"\n$insertBeforeIndex++\n\n$d.Paragraphs($insertBeforeIndex).Range.InsertParagraphBefore()\n$d.Paragraphs($insertBeforeIndex).Range.Text = \"RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; A/B Testing & Conversion Optimization\"\n$insertBeforeIndex++\n\n$d.Paragraphs($insertBeforeIndex).Range.InsertParagraphBefore()\n$d.Paragraphs($insertBeforeIndex).Range.Text = \"COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Client Relationship Management & Business Development\"\n"}
